$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Swap the match data (columns F..V) between row 31 (SC Ashdod vs Hapoel
# Hadera) and row 32 (Hapoel Haifa vs Maccabi Bnei Raina). The index/meta
# columns A..E stay attached to their original row.
# ---------------------------------------------------------------------------

$row31 = @{
    F = "Hapoel Haifa"
    G = 1
    H = "Maccabi Bnei Raina"
    I = 2
    J = 2.07
    K = "23/09/2023 18:13"
    L = 2.11
    M = "30/09/2023 18:44"
    N = 3.37
    O = "23/09/2023 18:13"
    P = 3.32
    Q = "30/09/2023 18:44"
    R = 3.63
    S = "23/09/2023 18:13"
    T = 3.72
    U = "30/09/2023 18:44"
    V = "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-haifa-maccabi-bnei-raina/xU1wcOQ1/"
}

$row32 = @{
    F = "SC Ashdod"
    G = 0
    H = "Hapoel Hadera"
    I = 1
    J = 2.13
    K = "23/09/2023 18:13"
    L = 1.81
    M = "30/09/2023 18:41"
    N = 3.34
    O = "23/09/2023 18:13"
    P = 3.58
    Q = "30/09/2023 18:44"
    R = 3.48
    S = "23/09/2023 18:13"
    T = 4.6
    U = "30/09/2023 18:44"
    V = "https://www.betexplorer.com/football/israel/ligat-ha-al/sc-ashdod-hapoel-hadera/jXkLL732/"
}

foreach ($col in $row31.Keys) {
    $ws.Range($col + "31").Value = $row31[$col]
}

foreach ($col in $row32.Keys) {
    $ws.Range($col + "32").Value = $row32[$col]
}

# ---------------------------------------------------------------------------
# Append a new row 37 with the Hapoel Tel Aviv vs Hapoel Haifa match.
# ---------------------------------------------------------------------------

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "israel"
$ws.Range("C37").Value = "ligat-ha-al"
$ws.Range("D37").Value = "2023-2024"
$ws.Range("E37").Value = 45262.58333333334
$ws.Range("F37").Value = "Hapoel Tel Aviv"
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = "Hapoel Haifa"
$ws.Range("I37").Value = 1
$ws.Range("J37").Value = 2.29
$ws.Range("K37").Value = "02/10/2023 19:42"
$ws.Range("L37").Value = 2.03
$ws.Range("M37").Value = "02/12/2023 13:53"
$ws.Range("N37").Value = 3.27
$ws.Range("O37").Value = "02/10/2023 19:42"
$ws.Range("P37").Value = 3.43
$ws.Range("Q37").Value = "02/12/2023 13:53"
$ws.Range("R37").Value = 3.21
$ws.Range("S37").Value = "02/10/2023 19:42"
$ws.Range("T37").Value = 3.84
$ws.Range("U37").Value = "02/12/2023 13:53"
$ws.Range("V37").Value = "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-tel-aviv-hapoel-haifa/lQ4QI5mL/"

# Match the formatting of column A (bold/centered index style) and column E
# (custom date/time number format) used by the rest of the sheet, by copying
# the formats down from the row above (reuses the existing style entries
# instead of creating new duplicate ones).
$ws.Range("A36").Copy() | Out-Null
$ws.Range("A37").PasteSpecial(-4122) | Out-Null

$ws.Range("E36").Copy() | Out-Null
$ws.Range("E37").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
